$d = $word.ActiveDocument

# Locate the paragraph containing the "LOQ4054" requirement text, then
# remove the three paragraphs that follow it (empty paragraph, the
# "Ver no Jupiter..." paragraph, and the "(c) 2020 ..." footer paragraph),
# while keeping the empty paragraph and page-break paragraph after them.

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOQ4054*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $delStart = $target.Range.End
    $next = $target.Next()
    # Skip over: empty paragraph, "Ver no Jupiter..." paragraph, "(c) 2020..." paragraph
    $next = $next.Next()
    $next = $next.Next()
    $delEnd = $next.Range.End

    $r = $d.Range($delStart, $delEnd)
    $r.Delete()
}
